$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 / J1 -------------------------------------------------
# Copy the header style (bold, bordered, centered) from the existing H1
# header cell so the new headers match the other header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data columns I2:I81 and J2:J81 --------------------------------------
$iValues = @(8,6,6,3,8,8,8,7,8,7,8,7,6,7,7,9,7,9,10,9,8,8,7,7,7,7,8,7,7,8,8,8,8,9,8,9,9,10,8,8,8,8,8,9,8,7,8,7,8,8,7,6,7,8,8,8,8,7,8,9,8,7,8,8,8,8,8,8,8,9,7,7,8,9,4,8,6,9,5,4)
$jValues = @(8,6,6,4,8,8,8,7,8,7,9,7,6,7,8,9,7,9,10,9,8,8,7,7,7,8,8,8,8,8,8,8,8,9,8,9,9,10,8,9,8,8,8,9,8,7,9,7,8,8,7,6,8,8,8,8,8,7,8,9,8,7,8,8,8,8,8,8,8,9,7,8,8,9,4,8,6,9,6,4)

$startRow = 2
for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $startRow + $n
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
